$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("C11").Value = -12.2813
$ws.Range("A12").Value = -21.3803
$ws.Range("C23").Value = -12.0628
$ws.Range("D24").Value = -8.671699999999989
$ws.Range("C28").Value = -13.89599999999999
$ws.Range("A32").Value = -21.421
$ws.Range("C32").Value = -12.7838
$ws.Range("C34").Value = -12.3444
$ws.Range("A36").Value = -20.0209
$ws.Range("A38").Value = -19.9686
$ws.Range("D38").Value = -7.678100000000001
$ws.Range("C42").Value = -12.8584
$ws.Range("A46").Value = -21.98440000000002
$ws.Range("D52").Value = -7.878900000000003
$ws.Range("A54").Value = -21.74409999999999
$ws.Range("C54").Value = -12.9889
$ws.Range("A55").Value = -22.25200000000001
$ws.Range("A67").Value = -21.56989999999996
$ws.Range("A69").Value = -21.69859999999997
$ws.Range("A72").Value = -21.77579999999999
$ws.Range("D78").Value = -7.505100000000001
$ws.Range("D83").Value = -9.081099999999999
$ws.Range("D85").Value = -8.738300000000002
$ws.Range("D86").Value = -8.6388
$ws.Range("A91").Value = -20.70139999999998
$ws.Range("D96").Value = -8.432199999999989
$ws.Range("C97").Value = -11.46
$ws.Range("A99").Value = -21.85679999999999
$ws.Range("C99").Value = -13.20489999999999
$ws.Range("C101").Value = -13.1018
$ws.Range("D103").Value = -7.820700000000003
$ws.Range("A104").Value = -21.38429999999999
